$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in I15 ("55.56%") first so it claims shared-string index 16 ---
# (the target workbook registers "55.56%" before "После ресайза" in
# sharedStrings.xml, so it must be the first new literal text written).
#
# I15 already carries a percentage number format (numFmtId 10). A plain
# Range.Value assignment of "55.56%" would be auto-parsed into the number
# 0.5556 (kept with that same percent format) instead of literal text, and
# forcing text via NumberFormat = "@" would mint a brand-new, unused style
# entry in styles.xml. To land literal text in I15 while keeping its
# existing style (s="3") untouched, build the text on a scratch cell via a
# formula that evaluates to a string, copy it, and paste only the *value*
# into I15 - the destination keeps its own formatting, only the stored
# value/type changes.
$scratch = $ws.Cells.Item(100, 26)
$scratch.Formula = '="55.56" & "%"'
$scratch.Copy() | Out-Null
$ws.Cells.Item(15, 9).PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null

# --- Fill in the rest of row 15 ---
$ws.Range("A15").Value = "После ресайза"
$ws.Range("B15").Value = "256x256"
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 450
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 12

# --- Move the active selection from H15 to A15 ---
$ws.Range("A15").Select() | Out-Null
